# Apply updated cryptocurrency price/volume data to sheet1 (ActiveSheet)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.573.49"
$ws.Range("E2").Value = "  -2.60%  "
$ws.Range("D3").Value = "1.753.63"
$ws.Range("E3").Value = "  -3.31%  "
$ws.Range("E4").Value = "  +0.12%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "324.15"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.42%  "
$ws.Range("E6").Value = "  +0.18%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.4451"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +1.81%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3621"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -1.15%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.07502"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -2.05%  "
$ws.Range("E10").Value = "  -5.53%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "1.107"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -2.84%  "
$ws.Range("E12").Value = "  +0.14%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "20.70"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -5.83%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "6.046"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -4.02%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "7.190"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -3.93%  "
$ws.Range("D16").Value = "1.753.13"
$ws.Range("E16").Value = "  -4.09%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "93.04"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -2.14%  "
$ws.Range("E18").Value = "  -1.40%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.06421"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -0.96%  "
$ws.Range("E20").Value = "  +0.17%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "17.06"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -1.81%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "5.842"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -6.37%  "
$ws.Range("D23").Value = "27.607.08"
$ws.Range("E23").Value = "  -2.49%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "11.26"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -2.43%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.100"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -0.31%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "162.86"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +0.95%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "20.43"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -1.43%  "
$ws.Range("D28").Value = "1.950.23"
$ws.Range("E28").Value = "  -4.00%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "2.130"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -6.35%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "125.74"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -2.45%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "1.087"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -9.72%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.09019"
$c.Style = "Normal"
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "3.639"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +1.63%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "5.551"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -7.53%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "12.11"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -6.32%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.02308"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -2.21%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.2095"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -3.49%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.6359"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -3.45%  "
$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "4.951"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -5.01%  "
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.05961"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -4.02%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "1.195"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +0.49%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "1.000"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +0.20%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "1.385"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -2.54%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "7.806"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -3.43%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "13.22"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -4.03%  "
$ws.Range("E46").Value = "  -0.53%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.5880"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -3.66%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "1.961"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -2.60%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "121.64"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -2.91%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "1.159"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +0.42%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.06860"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -1.86%  "
